# Actualización desde MV -datos-
# Append new daily rows (07-09-2021 .. 14-09-2021) to the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("07-09-2021", 904, 5143, 77),
    @("08-09-2021", 895, 5093, 76),
    @("09-09-2021", 884, 5031, 75),
    @("10-09-2021", 880, 5007, 75),
    @("11-09-2021", 880, 5007, 75),
    @("12-09-2021", 880, 5007, 75),
    @("13-09-2021", 881, 5015, 75),
    @("14-09-2021", 887, 5052, 75)
)

$startRow = 251
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    # Column A holds the date written as plain text (e.g. "07-09-2021"),
    # not a real date value. Entering that literal string directly would
    # let Excel auto-convert it into a date serial, so instead enter it as
    # a text-literal formula and immediately collapse it down to its
    # static value (Copy + PasteSpecial values-only) to land a plain
    # shared-string cell, exactly like the rest of the column.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Formula = '="' + $rowData[0] + '"'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)  # xlPasteValues

    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
}

$excel.CutCopyMode = $false
